# Update the "time_taken" column (F2:F128) on the "data" sheet with refreshed
# query timestamps.
$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

$timestamps = @(
    "2021-10-05 14:35:16.353422",
    "2021-10-05 14:35:16.353430",
    "2021-10-05 14:35:16.353433",
    "2021-10-05 14:35:16.353436",
    "2021-10-05 14:35:16.353438",
    "2021-10-05 14:35:16.353441",
    "2021-10-05 14:35:16.353444",
    "2021-10-05 14:35:16.353446",
    "2021-10-05 14:35:16.353449",
    "2021-10-05 14:35:16.353452",
    "2021-10-05 14:35:16.353454",
    "2021-10-05 14:35:16.353457",
    "2021-10-05 14:35:16.353459",
    "2021-10-05 14:35:16.353462",
    "2021-10-05 14:35:16.353464",
    "2021-10-05 14:35:16.353467",
    "2021-10-05 14:35:16.353469",
    "2021-10-05 14:35:16.353472",
    "2021-10-05 14:35:16.353475",
    "2021-10-05 14:35:16.353477",
    "2021-10-05 14:35:16.353480",
    "2021-10-05 14:35:16.353482",
    "2021-10-05 14:35:16.353485",
    "2021-10-05 14:35:16.353487",
    "2021-10-05 14:35:16.353490",
    "2021-10-05 14:35:16.353492",
    "2021-10-05 14:35:16.353495",
    "2021-10-05 14:35:16.353497",
    "2021-10-05 14:35:16.353500",
    "2021-10-05 14:35:16.353502",
    "2021-10-05 14:35:16.353505",
    "2021-10-05 14:35:16.353507",
    "2021-10-05 14:35:16.353510",
    "2021-10-05 14:35:16.353513",
    "2021-10-05 14:35:16.353515",
    "2021-10-05 14:35:16.353518",
    "2021-10-05 14:35:16.353520",
    "2021-10-05 14:35:16.353523",
    "2021-10-05 14:35:16.353525",
    "2021-10-05 14:35:16.353528",
    "2021-10-05 14:35:16.353531",
    "2021-10-05 14:35:16.353533",
    "2021-10-05 14:35:16.353536",
    "2021-10-05 14:35:16.353538",
    "2021-10-05 14:35:16.353541",
    "2021-10-05 14:35:16.353543",
    "2021-10-05 14:35:16.353546",
    "2021-10-05 14:35:16.353548",
    "2021-10-05 14:35:16.353551",
    "2021-10-05 14:35:16.353553",
    "2021-10-05 14:35:16.353556",
    "2021-10-05 14:35:16.353558",
    "2021-10-05 14:35:16.353561",
    "2021-10-05 14:35:16.353564",
    "2021-10-05 14:35:16.353566",
    "2021-10-05 14:35:16.353569",
    "2021-10-05 14:35:16.353571",
    "2021-10-05 14:35:16.353574",
    "2021-10-05 14:35:16.353576",
    "2021-10-05 14:35:16.353579",
    "2021-10-05 14:35:16.353581",
    "2021-10-05 14:35:16.353584",
    "2021-10-05 14:35:16.353586",
    "2021-10-05 14:35:16.353589",
    "2021-10-05 14:35:16.353593",
    "2021-10-05 14:35:16.353596",
    "2021-10-05 14:35:16.353598",
    "2021-10-05 14:35:16.353600",
    "2021-10-05 14:35:16.353603",
    "2021-10-05 14:35:16.353605",
    "2021-10-05 14:35:16.353608",
    "2021-10-05 14:35:16.353610",
    "2021-10-05 14:35:16.353613",
    "2021-10-05 14:35:16.353615",
    "2021-10-05 14:35:16.353618",
    "2021-10-05 14:35:16.353620",
    "2021-10-05 14:35:16.353625",
    "2021-10-05 14:35:16.353628",
    "2021-10-05 14:35:16.353630",
    "2021-10-05 14:35:16.353633",
    "2021-10-05 14:35:16.353635",
    "2021-10-05 14:35:16.353637",
    "2021-10-05 14:35:16.353640",
    "2021-10-05 14:35:16.353642",
    "2021-10-05 14:35:16.353645",
    "2021-10-05 14:35:16.353647",
    "2021-10-05 14:35:16.353650",
    "2021-10-05 14:35:16.353652",
    "2021-10-05 14:35:16.353655",
    "2021-10-05 14:35:16.353657",
    "2021-10-05 14:35:16.353660",
    "2021-10-05 14:35:16.353662",
    "2021-10-05 14:35:16.353666",
    "2021-10-05 14:35:16.353669",
    "2021-10-05 14:35:16.353672",
    "2021-10-05 14:35:16.353674",
    "2021-10-05 14:35:16.353677",
    "2021-10-05 14:35:16.353679",
    "2021-10-05 14:35:16.353681",
    "2021-10-05 14:35:16.353684",
    "2021-10-05 14:35:16.353686",
    "2021-10-05 14:35:16.353689",
    "2021-10-05 14:35:16.353691",
    "2021-10-05 14:35:16.353694",
    "2021-10-05 14:35:16.353696",
    "2021-10-05 14:35:16.353699",
    "2021-10-05 14:35:16.353701",
    "2021-10-05 14:35:16.353704",
    "2021-10-05 14:35:16.353708",
    "2021-10-05 14:35:16.353711",
    "2021-10-05 14:35:16.353714",
    "2021-10-05 14:35:16.353716",
    "2021-10-05 14:35:16.353719",
    "2021-10-05 14:35:16.353721",
    "2021-10-05 14:35:16.353724",
    "2021-10-05 14:35:16.353726",
    "2021-10-05 14:35:16.353729",
    "2021-10-05 14:35:16.353731",
    "2021-10-05 14:35:16.353734",
    "2021-10-05 14:35:16.353737",
    "2021-10-05 14:35:16.353739",
    "2021-10-05 14:35:16.353742",
    "2021-10-05 14:35:16.353744",
    "2021-10-05 14:35:16.353747",
    "2021-10-05 14:35:16.353749",
    "2021-10-05 14:35:16.353751",
    "2021-10-05 14:35:16.353754"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 6).Value = $timestamps[$i]
}

# Add a new "metadata" worksheet after "data", describing the panel query
# that produced the "data" sheet.
$meta = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSheet)
$meta.Name = "metadata"

# Header row (bold / bordered, matching the "data" sheet's header style).
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"
$meta.Range("B1:G1").Font.Bold = $true
$meta.Range("B1:G1").Borders.LineStyle = 1
$meta.Range("B1:G1").HorizontalAlignment = -4108
$meta.Range("B1:G1").VerticalAlignment = -4160

# Data row.
$meta.Range("A2").Value = 0
$meta.Range("A2").Font.Bold = $true
$meta.Range("A2").Borders.LineStyle = 1
$meta.Range("A2").HorizontalAlignment = -4108
$meta.Range("A2").VerticalAlignment = -4160

$meta.Range("B2").Value = "Polydactyly"
$meta.Range("C2").Value = 159
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "0.240"
$meta.Range("E2").Value = "2021-09-22T01:14:48.995373Z"
$meta.Range("F2").Value = "2021-10-05 14:35:16.349970"
$meta.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/159/?format=json"

# Leave the original "data" sheet as the active tab.
$dataSheet.Activate()
